$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook's config table had two DMS "folder/document url" rows
# (DMSReturnFolderUrl, DMSDocumentUrl) removed entirely, and the
# DMSEmailTitle / DMSExcelReturnTitle rows' Value text updated to the new
# "Waste Return ..." wording. Deleting the two rows (both originally at
# row 22) shifts everything below up by two, matching the table/dimension
# shrinking from A1:C29 to A1:C27.
$ws.Rows(22).Delete()
$ws.Rows(22).Delete()

# DMSExcelReturnTitle (now row 24) - Value
$ws.Cells.Item(24, 2).Value = "Waste Return {0} {1}"

# DMSEmailTitle (now row 22) - Value
$ws.Cells.Item(22, 2).Value = "Waste Return Correspondence {0} {1} - Email and Submission"

# Update the visible selection/scroll position to match the edited file:
# active cell A22, selected row A22:XFD22, scrolled so row 18 is at top.
$aw = $excel.ActiveWindow
$aw.ScrollRow = 18
$aw.ScrollColumn = 1
$ws.Range("A22:XFD22").Select()
